# Update "想去人数" (interested count, column F) and a couple of
# "最低票价" (min ticket price, column G) values across the three
# worksheets that carry this data: 展览, 演出, and 全部类型.
# Values correspond to a refreshed data pull (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibition) sheet ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 237
$ws1.Range("F3").Value  = 1060
$ws1.Range("F4").Value  = 529
$ws1.Range("F5").Value  = 13948
$ws1.Range("F7").Value  = 563
$ws1.Range("F8").Value  = 220
$ws1.Range("F9").Value  = 1797
$ws1.Range("F11").Value = 144
$ws1.Range("F12").Value = 95
$ws1.Range("F13").Value = 52
$ws1.Range("F14").Value = 541
$ws1.Range("F16").Value = 3
$ws1.Range("F18").Value = 14008
$ws1.Range("G18").Value = 60
$ws1.Range("F19").Value = 370
$ws1.Range("F20").Value = 633
$ws1.Range("F21").Value = 14991
$ws1.Range("F23").Value = 8300
$ws1.Range("F24").Value = 279
$ws1.Range("F26").Value = 29
$ws1.Range("F30").Value = 3
$ws1.Range("F31").Value = 12
$ws1.Range("F32").Value = 2
$ws1.Range("F34").Value = 1041
$ws1.Range("F35").Value = 22
$ws1.Range("F36").Value = 21
$ws1.Range("F37").Value = 23
$ws1.Range("F39").Value = 9
$ws1.Range("F43").Value = 394
$ws1.Range("F45").Value = 5107

# --- 演出 (Performance) sheet ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 50

# --- 全部类型 (All types) sheet ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 237
$ws4.Range("F3").Value  = 1060
$ws4.Range("F4").Value  = 529
$ws4.Range("F5").Value  = 13948
$ws4.Range("F7").Value  = 563
$ws4.Range("F8").Value  = 220
$ws4.Range("F9").Value  = 1797
$ws4.Range("F11").Value = 144
$ws4.Range("F12").Value = 95
$ws4.Range("F13").Value = 52
$ws4.Range("F14").Value = 541
$ws4.Range("F16").Value = 3
$ws4.Range("F18").Value = 14008
$ws4.Range("G18").Value = 60
$ws4.Range("F19").Value = 370
$ws4.Range("F20").Value = 633
$ws4.Range("F21").Value = 14991
$ws4.Range("F23").Value = 8300
$ws4.Range("F24").Value = 279
$ws4.Range("F26").Value = 29
$ws4.Range("F30").Value = 3
$ws4.Range("F31").Value = 12
$ws4.Range("F32").Value = 2
$ws4.Range("F34").Value = 1041
$ws4.Range("F35").Value = 22
$ws4.Range("F36").Value = 21
$ws4.Range("F37").Value = 23
$ws4.Range("F38").Value = 50
$ws4.Range("F41").Value = 9
$ws4.Range("F45").Value = 394
$ws4.Range("F47").Value = 5107

$wb.Save()
